$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell -> new value, matching the diff exactly.
$updates = @{
    "G2" = 13.86674266666667
    "H2" = 41.600228
    "I2" = 0.9894894203110381
    "J2" = 0.989489420311038
    "M2" = 5.079146
    "N2" = 15.237438
    "O2" = 0.707738273396773
    "P2" = 0.7077382733967729
    "Q2" = 70.43121054842933
    "R2" = 633.880894935864
    "S2" = 0.700299533875308
    "T2" = 0.7002995338753077

    "G3" = 13.86674266666667
    "H3" = 41.600228
    "I3" = 0.9894894203110381
    "J3" = 0.989489420311038
    "O3" = 0.2087923026002514
    "P3" = 0.2087923026002514
    "Q3" = 20.77815370186355
    "R3" = 187.003383316772
    "S3" = 0.2065977744653296
    "T3" = 0.2065977744653296

    "G4" = 13.86674266666667
    "H4" = 41.600228
    "I4" = 0.9894894203110381
    "J4" = 0.989489420311038
    "M4" = 0.5990256666666667
    "N4" = 1.797077
    "O4" = 0.08346942400297561
    "P4" = 0.08346942400297561
    "Q4" = 8.306534770395112
    "R4" = 74.758812933556
    "S4" = 0.08259211197040059
    "T4" = 0.08259211197040058

    "I5" = 0.006544770642485517
    "J5" = 0.006544770642485517
    "M5" = 5.079146
    "N5" = 15.237438
    "O5" = 0.707738273396773
    "P5" = 0.7077382733967729
    "Q5" = 0.4658524989253333
    "R5" = 4.192672490328
    "S5" = 0.004631984674290588
    "T5" = 0.004631984674290587

    "I6" = 0.006544770642485517
    "J6" = 0.006544770642485517
    "O6" = 0.2087923026002514
    "P6" = 0.2087923026002514
    "S6" = 0.001366497732435078
    "T6" = 0.001366497732435078

    "I7" = 0.006544770642485517
    "J7" = 0.006544770642485517
    "M7" = 0.5990256666666667
    "N7" = 1.797077
    "O7" = 0.08346942400297561
    "P7" = 0.08346942400297561
    "Q7" = 0.05494183544577778
    "R7" = 0.494476519012
    "S7" = 0.0005462882357598507
    "T7" = 0.0005462882357598507

    "G8" = 0.05557699999999999
    "H8" = 0.166731
    "I8" = 0.003965809046476372
    "J8" = 0.003965809046476372
    "M8" = 5.079146
    "N8" = 15.237438
    "O8" = 0.707738273396773
    "P8" = 0.7077382733967729
    "Q8" = 0.282283697242
    "R8" = 2.540553275178
    "S8" = 0.00280675484717449
    "T8" = 0.00280675484717449

    "G9" = 0.05557699999999999
    "H9" = 0.166731
    "I9" = 0.003965809046476372
    "J9" = 0.003965809046476372
    "O9" = 0.2087923026002514
    "P9" = 0.2087923026002514
    "Q9" = 0.08327748455766665
    "R9" = 0.7494973610189999
    "S9" = 0.0008280304024867092
    "T9" = 0.0008280304024867092

    "G10" = 0.05557699999999999
    "H10" = 0.166731
    "I10" = 0.003965809046476372
    "J10" = 0.003965809046476372
    "M10" = 0.5990256666666667
    "N10" = 1.797077
    "O10" = 0.08346942400297561
    "P10" = 0.08346942400297561
    "Q10" = 0.03329204947633333
    "R10" = 0.299628445287
    "S10" = 0.0003310237968151727
    "T10" = 0.0003310237968151727
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
